# Applies the "added requirement openpyxl" change:
#  - sheet1 ("test_file"): insert a new "sex" column between "age" and
#    "customer_type", populate it for every data row.
#  - sheet2 ("Sheet2"): trim the extra sample rows down to just the header
#    plus two data rows, and leave Sheet2 as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1: insert new column E ("sex") ------------------------------
$ws1.Columns.Item(5).Insert()

$ws1.Cells.Item(1, 5).Value = "sex"

$sexValues = @{
    2  = "Male"
    3  = "Other"
    4  = "Female"
    5  = "Male"
    6  = "Female"
    7  = "Male"
    8  = "Female"
    9  = "Male"
    10 = "Female"
    11 = "Male"
    12 = "Male"
    13 = "Male"
    14 = "Female"
    15 = "Female"
    16 = "Female"
    17 = "Female"
    18 = "Female"
    19 = "Male"
    20 = "Female"
    21 = "Male"
    22 = "Male"
    23 = "Female"
    24 = "Female"
    25 = "Female"
    26 = "Male"
    27 = "Male"
    28 = "Male"
    29 = "Male"
    30 = "Female"
}

foreach ($r in $sexValues.Keys) {
    $ws1.Cells.Item($r, 5).Value = $sexValues[$r]
}

# --- sheet2: drop the repeated sample rows (keep header + 2 rows) -----
$ws2.Rows("4:7").Delete()

# --- make Sheet2 the active tab / selection ---------------------------
$ws2.Activate()
$ws2.Range("A3").Select()
